# Insert a new data row at row 405 (shifting existing rows 405-458 down to 406-459)
# and populate it with a new "Berenjena" price record for Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 405; everything below shifts down one row.
$ws.Rows.Item(405).Insert()

# Populate the newly inserted row 405 with the new record's data.
$ws.Cells.Item(405, 1).Value2 = 3
$ws.Cells.Item(405, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(405, 3).Value = "Coquimbo"
$ws.Cells.Item(405, 4).Value2 = 45077
$ws.Cells.Item(405, 5).Value2 = 5
$ws.Cells.Item(405, 6).Value2 = 100112001
$ws.Cells.Item(405, 7).Value = "Berenjena"
$ws.Cells.Item(405, 8).Value = "Sin especificar"
$ws.Cells.Item(405, 9).Value = "Primera"
$ws.Cells.Item(405, 10).Value2 = 85
$ws.Cells.Item(405, 11).Value2 = 7500
$ws.Cells.Item(405, 12).Value2 = 8000
$ws.Cells.Item(405, 13).Value2 = 7765
$ws.Cells.Item(405, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(405, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(405, 16).Value2 = 129
$ws.Cells.Item(405, 17).Value2 = 60
$ws.Cells.Item(405, 18).Value = "Hortaliza"

# Apply the same number format (date) used by the other rows in column D.
$ws.Cells.Item(405, 4).NumberFormat = $ws.Cells.Item(406, 4).NumberFormat
